# Update "想去人数" (F) and, in two spots, "最低票价" (G) values across the
# 展览 (Exhibition), 演出 (Performance) and 全部类型 (All types) sheets to
# match the freshly re-scraped numbers.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "展览" (Exhibition) - column F only, rows 3-17
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$exhibitionF = @{
    3  = 12137
    4  = 51
    5  = 234
    6  = 373
    7  = 230
    8  = 12042
    9  = 505
    10 = 1186
    11 = 113
    12 = 597
    13 = 2805
    14 = 5941
    15 = 134
    16 = 3560
    17 = 205
}
foreach ($row in $exhibitionF.Keys) {
    $ws1.Cells.Item($row, 6).Value = $exhibitionF[$row]
}

# ---------------------------------------------------------------------
# Sheet "演出" (Performance)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(2, 6).Value = 578
$ws2.Cells.Item(2, 7).Value = "不可售"
$ws2.Cells.Item(5, 6).Value = 5

# ---------------------------------------------------------------------
# Sheet "全部类型" (All types) - union of the two sheets above
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(2, 6).Value = 578
$ws4.Cells.Item(2, 7).Value = "不可售"

$allTypesF = @{
    5  = 12137
    6  = 51
    7  = 234
    9  = 373
    10 = 230
    11 = 12042
    12 = 505
    13 = 1186
    14 = 113
    15 = 597
    16 = 2805
    17 = 5
    18 = 5941
    19 = 134
    20 = 3560
    21 = 205
}
foreach ($row in $allTypesF.Keys) {
    $ws4.Cells.Item($row, 6).Value = $allTypesF[$row]
}
